$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new day's data to row 32 (date serial 45981 = 2025-11-20)
$ws.Range("A32").Value = 45981
$ws.Range("B32").Value = 597
$ws.Range("C32").Value = 20
$ws.Range("D32").Value = 577

# Move the active selection down to the newly-added row, matching the
# author's last selection state (A32:D32, active cell A32)
$ws.Range("A32:D32").Select()

Write-Output "done"
